$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 7 (keep only header + one data row)
$ws.Range("A3:D7").EntireRow.Delete()

# Update row 2 values with new gene/probe data
$ws.Range("A2").Value = "cg05072008"
$ws.Range("B2").Value = "FIGNL1"
$ws.Range("C2").Value = "auto"
$ws.Range("D2").Value = "auto"

# Update active selection to B2 as in target file
$ws.Range("B2").Select()
